# edit.ps1 - applies the target edits to the presentation via PowerPoint COM automation.
#
# Summary of changes:
#  1. Insert a brand-new "Title and Content" slide at position 4 (pushing the
#     existing "Plans" / "So far" / "Summary" / "Bibliography" slides down by
#     one), titled "Feasibility studies" with three bullet paragraphs
#     (separated by blank paragraphs).
#  2. On the "Aims" slide, fix the "electromagatism" typo -> "electromagnetism"
#     and tidy up the following run's leading space.
#  3. On the "Introduction" slide, rename the title to "Theory" and add a
#     "Bloch sphere" line (plus two blank paragraphs) to the body placeholder.

$p = $ppt.ActivePresentation

# --- 1. Insert the new "Feasibility studies" slide at index 4 -------------
# Layout 2 == ppLayoutText ("Title, Content"), matching the other
# title+bullets slides in this deck (e.g. "Plans").
$newSlide = $p.Slides.Add(4, 2)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Feasibility studies"

$bodyTr = $newSlide.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "Driving the rotor of the gyroscope with a current`r`rWith small conducting magnets`r`rTachometer and motor arrangement"

# --- 2. Fix "Aims" slide wording -------------------------------------------
$aims = $p.Slides.Item(2)
$aimsBody = $aims.Shapes.Item(2).TextFrame.TextRange
$aimsPara2 = $aimsBody.Paragraphs(2, 1)
$aimsPara2.Runs(2, 1).Text = "electromagnetism "
$aimsPara2.Runs(3, 1).Text = "causing the gyroscope to process"

# --- 3. Rename "Introduction" -> "Theory" and add "Bloch sphere" ----------
$theory = $p.Slides.Item(3)
$theory.Shapes.Item(1).TextFrame.TextRange.Text = "Theory"
$theory.Shapes.Item(2).TextFrame.TextRange.Text = "Bloch sphere`r`r"
